$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("boroughexpedition.com", "401"),
    @("communitynavigate.com", "401"),
    @("communityprobe.com", "401"),
    @("communitypilots.com", "401"),
    @("districtexplorer.com", "401"),
    @("localeguided.com", "200"),
    @("communityroamers.com", "401"),
    @("communitytrackers.com", "401"),
    @("localequest.xyz", "500"),
    @("districtsearchers.com", "401"),
    @("localequest.com", "200"),
    @("metropathfinders.com", "401"),
    @("regionaldetective.com", "200"),
    @("regionalprobes.com", "200"),
    @("regionalscavenger.com", "200"),
    @("townrovers.com", "401"),
    @("townsleuth.com", "401"),
    @("vicinityfinders.com", "401")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Formula = "'" + $data[$i][1]
    $ws.Cells.Item($row, 2).Style = "Normal"
}
